$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '91.311.45'
$ws.Range("E2").Value = '  +0.98%  '

# Row 3
$ws.Range("D3").Value = '3.162.34'
$ws.Range("E3").Value = '  +2.30%  '

# Row 4
$ws.Range("E4").Value = '  +0.33%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.05'
$ws.Range("E5").Value = '  +2.16%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.11'
$ws.Range("E6").Value = '  -0.25%  '

# Row 7
$ws.Range("E7").Value = '  +4.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.369'
$ws.Range("E8").Value = '  +0.67%  '

# Row 9
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("D10").Value = '3.164.42'
$ws.Range("E10").Value = '  +2.51%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.740'
$ws.Range("E11").Value = '  +1.14%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.204'
$ws.Range("E12").Value = '  +3.61%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000248'
$ws.Range("E13").Value = '  -1.44%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.43'
$ws.Range("E14").Value = '  -2.08%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.49'
$ws.Range("E15").Value = '  +0.27%  '

# Row 16
$ws.Range("D16").Value = '91.277.91'
$ws.Range("E16").Value = '  +1.37%  '

# Row 17
$ws.Range("D17").Value = '3.734.28'
$ws.Range("E17").Value = '  +2.00%  '

# Row 18
$ws.Range("D18").Value = '3.221.48'
$ws.Range("E18").Value = '  +4.86%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.72'
$ws.Range("E19").Value = '  -3.99%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.51'
$ws.Range("E20").Value = '  +10.71%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000209'
$ws.Range("E21").Value = '  -3.63%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.82'
$ws.Range("E22").Value = '  +4.36%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '443.31'
$ws.Range("E23").Value = '  +1.67%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.21'
$ws.Range("E24").Value = '  +3.34%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.81'
$ws.Range("E25").Value = '  -2.14%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.17'
$ws.Range("E26").Value = '  -0.06%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.02'
$ws.Range("E27").Value = '  -0.61%  '

# Row 28
$ws.Range("D28").Value = '3.271.38'
$ws.Range("E28").Value = '  +0.68%  '

# Row 29
$ws.Range("E29").Value = '  -0.02%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.136'
$ws.Range("E30").Value = '  +52.40%  '

# Row 31
$ws.Range("E31").Value = '  +18.43%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.172'
$ws.Range("E32").Value = '  +8.27%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '9.30'
$ws.Range("E33").Value = '  -0.99%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.165'
$ws.Range("E34").Value = '  +7.56%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.86'
$ws.Range("E35").Value = '  +9.57%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '26.51'
$ws.Range("E36").Value = '  +2.86%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.16'
$ws.Range("E37").Value = '  +21.11%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '511.45'
$ws.Range("E38").Value = '  +1.59%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.35'
$ws.Range("E39").Value = '  +5.06%  '

# Row 40
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.459'
$ws.Range("E40").Value = '  +13.37%  '

# Row 41
$ws.Range("B41").Value = 'PancakeSwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.92'
$ws.Range("E41").Value = '  +0.55%  '

# Row 42
$ws.Range("B42").Value = 'Binance-PegBSC-USD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.813'
$ws.Range("E42").Value = '  -18.44%  '

# Row 43
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.46'
$ws.Range("E43").Value = '  -9.62%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.16'
$ws.Range("E44").Value = '  +0.01%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.717'
$ws.Range("E46").Value = '  +3.71%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.92'
$ws.Range("E47").Value = '  +0.85%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '155.86'
$ws.Range("E48").Value = '  +2.12%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.39'
$ws.Range("E49").Value = '  +3.48%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.45'
$ws.Range("E50").Value = '  +0.73%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0327'
$ws.Range("E51").Value = '  +13.07%  '
